$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hyperlink's internal location target (edit in place via
# enumeration - Hyperlinks.Item(n) does not resolve correctly here)
foreach ($h in $ws.Hyperlinks) {
    $h.SubAddress = "result-for-submission"
}

# Update the hyperlink display text (shared string)
$ws.Range("H3").Value = "results"

# Update the H2 value
$ws.Range("H2").Value = 5813768

# Update the selected cell (active cell) to F20
$ws.Range("F20").Select()
